$d = $word.ActiveDocument
$null = $d.Content.Find.Execute("657×3=1971", $true, $false, $false, $false, $false, $true, 1, $false, "892×7=6244", 2)
$null = $d.Content.Find.Execute("227×9=2043", $true, $false, $false, $false, $false, $true, 1, $false, "319×2=638", 2)
$null = $d.Content.Find.Execute("116×7=812", $true, $false, $false, $false, $false, $true, 1, $false, "184×7=1288", 2)
$null = $d.Content.Find.Execute("782×9=7038", $true, $false, $false, $false, $false, $true, 1, $false, "189×2=378", 2)
$null = $d.Content.Find.Execute("852×4=3408", $true, $false, $false, $false, $false, $true, 1, $false, "321×6=1926", 2)
$null = $d.Content.Find.Execute("328×9=2952", $true, $false, $false, $false, $false, $true, 1, $false, "978×6=5868", 2)
$null = $d.Content.Find.Execute("591×7=4137", $true, $false, $false, $false, $false, $true, 1, $false, "107×6=642", 2)
$null = $d.Content.Find.Execute("153×4=612", $true, $false, $false, $false, $false, $true, 1, $false, "866×7=6062", 2)
$null = $d.Content.Find.Execute("529×2=1058", $true, $false, $false, $false, $false, $true, 1, $false, "686×6=4116", 2)
$null = $d.Content.Find.Execute("846×4=3384", $true, $false, $false, $false, $false, $true, 1, $false, "264×9=2376", 2)
$null = $d.Content.Find.Execute("149×6=894", $true, $false, $false, $false, $false, $true, 1, $false, "162×3=486", 2)
$null = $d.Content.Find.Execute("196×4=784", $true, $false, $false, $false, $false, $true, 1, $false, "472×5=2360", 2)
$null = $d.Content.Find.Execute("522×7=3654", $true, $false, $false, $false, $false, $true, 1, $false, "340×5=1700", 2)
$null = $d.Content.Find.Execute("217×8=1736", $true, $false, $false, $false, $false, $true, 1, $false, "101×9=909", 2)
$null = $d.Content.Find.Execute("286×7=2002", $true, $false, $false, $false, $false, $true, 1, $false, "114×8=912", 2)
$null = $d.Content.Find.Execute("536×6=3216", $true, $false, $false, $false, $false, $true, 1, $false, "793×7=5551", 2)
$null = $d.Content.Find.Execute("906×8=7248", $true, $false, $false, $false, $false, $true, 1, $false, "161×3=483", 2)
$null = $d.Content.Find.Execute("906×6=5436", $true, $false, $false, $false, $false, $true, 1, $false, "520×2=1040", 2)
$null = $d.Content.Find.Execute("599×2=1198", $true, $false, $false, $false, $false, $true, 1, $false, "684×2=1368", 2)
$null = $d.Content.Find.Execute("493×7=3451", $true, $false, $false, $false, $false, $true, 1, $false, "593×9=5337", 2)
$null = $d.Content.Find.Execute("117×6=702", $true, $false, $false, $false, $false, $true, 1, $false, "296×5=1480", 2)
$null = $d.Content.Find.Execute("870×2=1740", $true, $false, $false, $false, $false, $true, 1, $false, "859×2=1718", 2)
$null = $d.Content.Find.Execute("247×2=494", $true, $false, $false, $false, $false, $true, 1, $false, "172×3=516", 2)
$null = $d.Content.Find.Execute("976×5=4880", $true, $false, $false, $false, $false, $true, 1, $false, "566×3=1698", 2)
$null = $d.Content.Find.Execute("714×4=2856", $true, $false, $false, $false, $false, $true, 1, $false, "243×4=972", 2)
